# Update the "PC build" spreadsheet:
#  - Parts sheet: swap out the old monitor/CPU/motherboard/RAM/GPU parts for a
#    new AMD Ryzen based build, and drop the Monitor line entirely.
#  - Category sheet: rename "Processors" -> "Processor" and drop the unused
#    "DVD" category.
#  - Webshop sheet: drop the unused "hobbycraft" / "hm" webshops.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Parts sheet
# ---------------------------------------------------------------------------
$parts = $wb.Worksheets.Item("Parts")

# Remove the "24 Inch Full HD Monitor" / "Monitor" row entirely.
$parts.Rows(2).Delete()

# Overwrite the remaining rows with the new build's parts/categories.
$parts.Range("A2").Value = "AMD Ryzen 5 7600X"
$parts.Range("A4").Value = "32GB DDR5 5600MHz"
$parts.Range("A6").Value = "NVIDIA GeForce RTX 4060"
$parts.Range("A3").Value = "Gigabyte B650 AORUS ELITE AX AMD B650 Chipset"
$parts.Range("B2").Value = "Processor"

$parts.Range("B3").Value = "Motherboard"
$parts.Range("B4").Value = "Memory "
$parts.Range("A5").Value = "ATX Mid Tower Gaming Case"
$parts.Range("B5").Value = "Computer Case"
$parts.Range("B6").Value = "Graphics Card"
$parts.Range("A7").Value = "1TB HDD"
$parts.Range("B7").Value = "hard drive "

# The long motherboard name wraps, so give A3 a wrap-text style.
$parts.Range("A3").WrapText = $true

# Widen column A to fit the longer part names.
$parts.Columns("A").ColumnWidth = 42.833333333333336

# ---------------------------------------------------------------------------
# Category sheet
# ---------------------------------------------------------------------------
$category = $wb.Worksheets.Item("Category")

# Drop the unused "DVD" category.
$category.Rows(5).Delete()

# Rename "Processors" -> "Processor".
$category.Range("A4").Value = "Processor"

$category.Rows(5).Select()

# ---------------------------------------------------------------------------
# Webshop sheet
# ---------------------------------------------------------------------------
$webshop = $wb.Worksheets.Item("Webshop")

# Drop the unused "hobbycraft" and "hm" webshops (rows 5 and 6).
$webshop.Range("A5:A6").EntireRow.Delete()

$webshop.Range("A5:A6").EntireRow.Select()

# ---------------------------------------------------------------------------
# Restore the Parts sheet as the active tab/selection, matching the saved
# view state (it was the tab visible when the workbook was last saved).
# ---------------------------------------------------------------------------
$parts.Activate()
$parts.Range("B3").Select()
